$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing credentials row (A2/B2) with new test data.
$ws.Range("A2").Value = "mcAngular"
$ws.Range("B2").Value = "3d6g4f7j5g8k"

# Add a new data row (A3/B3) for data-driven login testing.
$ws.Range("A3").Value = "tmtmoney"
$ws.Range("B3").Value = "oy06ri94uw73"

# Match the formatting of the header/credential rows above for the new row.
$ws.Range("A2:B2").Copy()
$ws.Range("A3:B3").PasteSpecial(-4122)

